# Auto commit - 09031537
# Adds a new maintenance-report row (row 27) to the 'Report' sheet, fixes
# wrap-text formatting that was missing on row 26's P/AC cells, extends the
# print area to include the new row, and updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Report")

# ---------------------------------------------------------------------
# 1) Fix existing row 26: P26 and AC26 should have wrap text enabled,
#    matching the style used by every other data row's P/AC columns.
# ---------------------------------------------------------------------
$ws.Range("P26").WrapText = $true
$ws.Range("AC26").WrapText = $true

# ---------------------------------------------------------------------
# 2) Insert a new data row (row 27) describing the latest work order.
#    First clone the formatting of the row above (row 25), which already
#    carries the correct style pattern for a freshly-added record, then
#    populate the cell values/text for row 27.
# ---------------------------------------------------------------------
$srcRow = $ws.Range("A25:AK25")
$dstRow = $ws.Range("A27:AK27")
$srcRow.Copy()
$dstRow.PasteSpecial(-4122)  # xlPasteFormats

$r = 27
$ws.Cells.Item($r, 1).Value = 25                          # A - 項次
$ws.Cells.Item($r, 2).Value = "服務"                       # B - 工作類別
$ws.Cells.Item($r, 3).Value = 2025090650                  # C - 台芝工作案號
$ws.Cells.Item($r, 6).Value = 4298                        # F - 門店編號
$ws.Cells.Item($r, 7).Value = "淡水後洲店"                  # G - 門店名稱
$ws.Cells.Item($r, 8).Value = "新北市淡水區"                # H - 縣市鄉鎮
$ws.Cells.Item($r, 17).Value = "THILF04298"               # Q - 設備號碼
$ws.Cells.Item($r, 18).Value = "新北一"                     # R - 負責部門
$ws.Cells.Item($r, 19).Value = "吳宗鴻"                     # S - 服務人員
$ws.Cells.Item($r, 20).Value = 1                          # T - 到場次數
$ws.Cells.Item($r, 21).Value = "已完工"                     # U - 完工結果
$ws.Cells.Item($r, 22).Value = "2025-09-03 15:31:31"      # V - 派修時間
$ws.Cells.Item($r, 23).Value = "2025-09-03 15:10:00"      # W - 到場時間
$ws.Cells.Item($r, 24).Value = "2025-09-03 15:30:00"      # X - 離場時間
$ws.Cells.Item($r, 26).Value = 0.3                        # Z - 處理工時
$ws.Cells.Item($r, 28).Value = "到場處理"                   # AB - 結案類別
$ws.Cells.Item($r, 29).Value = "PMQ3"                     # AC - 工作內容
$ws.Cells.Item($r, 30).Value = "O"                        # AD - 保養
$ws.Cells.Item($r, 37).Value = "O"                        # AK - 客戶簽章

# ---------------------------------------------------------------------
# 3) Extend the print area to cover the newly added row.
# ---------------------------------------------------------------------
$ws.PageSetup.PrintArea = '$A$1:$AK$27'

# ---------------------------------------------------------------------
# 4) Update the active selection to the new row's first cell.
# ---------------------------------------------------------------------
[void]$ws.Range("A27").Select()
